$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$prefix = '<jt:rickroll value="'
$label  = 'Additional Help'
$suffix = '"/>'

$ws.Range("A3").Value = "$prefix$label$suffix"

$startPos = $prefix.Length + 1
$chars = $ws.Range("A3").Characters($startPos, $label.Length)
$chars.Font.Underline = $true
$chars.Font.Color = 16711680

$suffixStart = $startPos + $label.Length
$suffixChars = $ws.Range("A3").Characters($suffixStart, $suffix.Length)
$suffixChars.Font.Color = 0

Write-Output "Added jt:rickroll hyperlink text to A3"
